# Regen save_data to use K instead of Strike#: update column G (K) values
# on the active worksheet for rows 2-12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 4
    3  = 5
    4  = 5
    5  = 6
    6  = 5
    7  = 9
    8  = 6
    9  = 8
    10 = 4
    11 = 2
    12 = 6
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
